$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 587, shifting existing rows 587:642 down to 588:643.
$ws.Rows.Item(587).Insert()

# Populate the newly inserted row 587 with the new data record.
$ws.Range("A587").Value = 3
$ws.Range("B587").Value = "Femacal de La Calera"
$ws.Range("C587").Value = "Coquimbo"
$ws.Range("D587").Value = 45132
$ws.Range("E587").Value = 5
$ws.Range("F587").Value = 100112031
$ws.Range("G587").Value = "Poroto verde"
$ws.Range("H587").Value = "Magnum"
$ws.Range("I587").Value = "Primera"
$ws.Range("J587").Value = 40
$ws.Range("K587").Value = 27000
$ws.Range("L587").Value = 27000
$ws.Range("M587").Value = 27000
$ws.Range("N587").Value = "$/malla 25 kilos"
$ws.Range("O587").Value = "Región de Arica y Parinacota"
$ws.Range("P587").Value = 1080
$ws.Range("Q587").Value = 25
$ws.Range("R587").Value = "Hortaliza"
